$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get redistributed between rows in each group.
$cols = @("A","B","D","E","F","G","H","Q","R","AC")

function Get-RowValues($row, $columns) {
    $vals = @{}
    foreach ($c in $columns) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($row, $columns, $vals) {
    foreach ($c in $columns) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

# --- Rows 100 and 101: swap content (including Observatorer / AX) ---
$colsWithAX = $cols + @("AX")
$row100 = Get-RowValues 100 $colsWithAX
$row101 = Get-RowValues 101 $colsWithAX

Set-RowValues 100 $colsWithAX $row101
Set-RowValues 101 $colsWithAX $row100

# --- Rows 105, 106, 107: rotate content (105<-106, 106<-107, 107<-105) ---
$row105 = Get-RowValues 105 $cols
$row106 = Get-RowValues 106 $cols
$row107 = Get-RowValues 107 $cols

Set-RowValues 105 $cols $row106
Set-RowValues 106 $cols $row107
Set-RowValues 107 $cols $row105

# --- Rows 112, 113, 114: rotate content (112<-114, 113<-112, 114<-113) ---
$row112 = Get-RowValues 112 $cols
$row113 = Get-RowValues 113 $cols
$row114 = Get-RowValues 114 $cols

Set-RowValues 112 $cols $row114
Set-RowValues 113 $cols $row112
Set-RowValues 114 $cols $row113
